# update Tue 2025-01-07 T-590 HF
# Fill in the "Pozn." (note) column D on the Results sheet with the
# proofreader initials for the rows that were checked, and move the
# sheet's scroll/selection state to reflect where work left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")
$ws.Activate()

# Rows that already had a reviewer note of "HAN" now also got checked by
# "Hnovot", so the note becomes "HAN, Hnovot".
$hanPlusHnovot = @(15, 20, 28, 32)
foreach ($r in $hanPlusHnovot) {
    $ws.Cells.Item($r, 4).Value = "HAN, Hnovot"
}

# Rows newly reviewed solely by "Hnovot".
$hnovotOnly = @(9, 25, 52, 67)
foreach ($r in $hnovotOnly) {
    $ws.Cells.Item($r, 4).Value = "HNovot"
}

# Rows newly reviewed by "PM" together with "Hnovot".
$pmPlusHnovot = @(40, 48)
foreach ($r in $pmPlusHnovot) {
    $ws.Cells.Item($r, 4).Value = "PM, Hnovot"
}

# Rows newly reviewed solely by "PM".
$pmOnly = @(50, 53, 57, 72, 76, 80, 82, 86)
foreach ($r in $pmOnly) {
    $ws.Cells.Item($r, 4).Value = "PM"
}

# Move the frozen-pane scroll position / active selection to where the
# reviewer was last working.
$ws.Range("A50").Select() | Out-Null
$ws.Range("C67").Select() | Out-Null
